$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was reported for this market/product; insert a new
# row above row 13 (shifting the existing rows 13-48 down to 14-49) and
# populate the new row with the latest observation.
$ws.Rows("13:13").Insert()

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 45044
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100101
$ws.Range("H13").Value = "Berries"
$ws.Range("I13").Value = 100101001
$ws.Range("J13").Value = "Arándano (blue)"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 12500
$ws.Range("Q13").Value = "$/bandeja 2 kilos"
$ws.Range("R13").Value = "Provincia de Curicó"
$ws.Range("S13").Value = 6250
$ws.Range("T13").Value = 2
